$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Nov 05 14:10:02 EST 2025"
$ws.Range("B3").Value = "Wed Nov 05 14:10:13 EST 2025"
$ws.Range("B4").Value = "Wed Nov 05 14:10:22 EST 2025"
$ws.Range("B5").Value = "Wed Nov 05 14:10:31 EST 2025"
$ws.Range("B6").Value = "Wed Nov 05 14:10:40 EST 2025"
$ws.Range("B7").Value = "Wed Nov 05 14:10:49 EST 2025"
$ws.Range("B8").Value = "Wed Nov 05 14:10:57 EST 2025"
$ws.Range("B9").Value = "Wed Nov 05 14:11:05 EST 2025"
$ws.Range("B10").Value = "Wed Nov 05 14:11:13 EST 2025"
$ws.Range("B11").Value = "Wed Nov 05 14:11:21 EST 2025"
$ws.Range("B12").Value = "Wed Nov 05 14:11:29 EST 2025"
$ws.Range("B13").Value = "Wed Nov 05 14:11:37 EST 2025"
$ws.Range("B14").Value = "Wed Nov 05 14:11:45 EST 2025"
$ws.Range("B15").Value = "Wed Nov 05 14:11:54 EST 2025"
$ws.Range("B16").Value = "Wed Nov 05 14:12:03 EST 2025"
$ws.Range("B17").Value = "Wed Nov 05 14:12:12 EST 2025"
$ws.Range("B18").Value = "Wed Nov 05 14:12:20 EST 2025"
$ws.Range("B19").Value = "Wed Nov 05 14:12:28 EST 2025"
$ws.Range("B20").Value = "Wed Nov 05 14:12:36 EST 2025"
$ws.Range("B21").Value = "Wed Nov 05 14:12:44 EST 2025"
$ws.Range("B22").Value = "Wed Nov 05 14:12:54 EST 2025"
$ws.Range("B23").Value = "Wed Nov 05 14:13:02 EST 2025"
$ws.Range("B24").Value = "Wed Nov 05 14:13:10 EST 2025"
$ws.Range("B25").Value = "Wed Nov 05 14:13:18 EST 2025"
$ws.Range("B26").Value = "Wed Nov 05 14:13:26 EST 2025"
$ws.Range("B27").Value = "Wed Nov 05 14:13:34 EST 2025"
$ws.Range("B28").Value = "Wed Nov 05 14:13:43 EST 2025"
$ws.Range("B29").Value = "Wed Nov 05 14:13:51 EST 2025"
$ws.Range("B30").Value = "Wed Nov 05 14:14:00 EST 2025"
$ws.Range("B31").Value = "Wed Nov 05 14:14:08 EST 2025"
$ws.Range("B32").Value = "Wed Nov 05 14:14:16 EST 2025"
$ws.Range("B33").Value = "Wed Nov 05 14:14:24 EST 2025"
$ws.Range("B34").Value = "Wed Nov 05 14:14:32 EST 2025"
$ws.Range("B36").Value = "Wed Nov 05 14:14:41 EST 2025"
$ws.Range("B37").Value = "Wed Nov 05 14:14:49 EST 2025"
$ws.Range("B38").Value = "Wed Nov 05 14:14:57 EST 2025"
$ws.Range("B39").Value = "Wed Nov 05 14:15:05 EST 2025"
$ws.Range("B40").Value = "Wed Nov 05 14:15:13 EST 2025"
$ws.Range("B42").Value = "Wed Nov 05 14:15:21 EST 2025"
$ws.Range("B43").Value = "Wed Nov 05 14:15:30 EST 2025"
$ws.Range("B44").Value = "Wed Nov 05 14:15:38 EST 2025"
$ws.Range("B45").Value = "Wed Nov 05 14:15:47 EST 2025"
$ws.Range("B47").Value = "Wed Nov 05 14:15:55 EST 2025"
$ws.Range("B48").Value = "Wed Nov 05 14:16:03 EST 2025"
$ws.Range("B49").Value = "Wed Nov 05 14:16:12 EST 2025"
$ws.Range("B50").Value = "Wed Nov 05 14:16:20 EST 2025"
$ws.Range("B51").Value = "Wed Nov 05 14:16:28 EST 2025"
$ws.Range("B52").Value = "Wed Nov 05 14:16:36 EST 2025"
$ws.Range("B53").Value = "Wed Nov 05 14:16:45 EST 2025"
$ws.Range("B54").Value = "Wed Nov 05 14:16:52 EST 2025"
$ws.Range("B55").Value = "Wed Nov 05 14:17:02 EST 2025"
$ws.Range("B56").Value = "Wed Nov 05 14:17:10 EST 2025"
$ws.Range("B57").Value = "Wed Nov 05 14:17:19 EST 2025"
$ws.Range("B59").Value = "Wed Nov 05 14:17:27 EST 2025"
